$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 16.37389066666667
$ws.Range("H2").Value = 49.121672
$ws.Range("I2").Value = 0.09466313117816218
$ws.Range("J2").Value = 0.09466313117816218
$ws.Range("M2").Value = 0.01339666666666667
$ws.Range("N2").Value = 0.04019
$ws.Range("O2").Value = 0.08393217762128816
$ws.Range("P2").Value = 0.08393217762128814
$ws.Range("Q2").Value = 0.2193555552977778
$ws.Range("R2").Value = 1.97419999768
$ws.Range("S2").Value = 0.007945282740232809
$ws.Range("T2").Value = 0.007945282740232807

# Row 3
$ws.Range("G3").Value = 16.37389066666667
$ws.Range("H3").Value = 49.121672
$ws.Range("I3").Value = 0.09466313117816218
$ws.Range("J3").Value = 0.09466313117816218
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.1260863333333333
$ws.Range("N3").Value = 0.378259
$ws.Range("O3").Value = 0.7899502755623498
$ws.Range("P3").Value = 0.7899502755623498
$ws.Range("Q3").Value = 2.064523836560889
$ws.Range("R3").Value = 18.580714529048
$ws.Range("S3").Value = 0.07477916655978407
$ws.Range("T3").Value = 0.07477916655978407

# Row 4
$ws.Range("G4").Value = 16.37389066666667
$ws.Range("H4").Value = 49.121672
$ws.Range("I4").Value = 0.09466313117816218
$ws.Range("J4").Value = 0.09466313117816218
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.02013
$ws.Range("N4").Value = 0.06039
$ws.Range("O4").Value = 0.1261175468163621
$ws.Range("P4").Value = 0.1261175468163621
$ws.Range("Q4").Value = 0.32960641912
$ws.Range("R4").Value = 2.96645777208
$ws.Range("S4").Value = 0.01193868187814529
$ws.Range("T4").Value = 0.01193868187814529

# Row 5
$ws.Range("I5").Value = 0.4193879037829277
$ws.Range("J5").Value = 0.4193879037829278
$ws.Range("M5").Value = 0.01339666666666667
$ws.Range("N5").Value = 0.04019
$ws.Range("O5").Value = 0.08393217762128816
$ws.Range("P5").Value = 0.08393217762128814
$ws.Range("Q5").Value = 0.9718151657833334
$ws.Range("R5").Value = 8.74633649205
$ws.Range("S5").Value = 0.03520014003252839
$ws.Range("T5").Value = 0.03520014003252839

# Row 6
$ws.Range("I6").Value = 0.4193879037829277
$ws.Range("J6").Value = 0.4193879037829278
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.1260863333333333
$ws.Range("N6").Value = 0.378259
$ws.Range("O6").Value = 0.7899502755623498
$ws.Range("P6").Value = 0.7899502755623498
$ws.Range("Q6").Value = 9.146499945111668
$ws.Range("R6").Value = 82.31849950600501
$ws.Range("S6").Value = 0.33129559016084
$ws.Range("T6").Value = 0.33129559016084

# Row 7
$ws.Range("I7").Value = 0.4193879037829277
$ws.Range("J7").Value = 0.4193879037829278
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.02013
$ws.Range("N7").Value = 0.06039
$ws.Range("O7").Value = 0.1261175468163621
$ws.Range("P7").Value = 0.1261175468163621
$ws.Range("Q7").Value = 1.46026170345
$ws.Range("R7").Value = 13.14235533105
$ws.Range("S7").Value = 0.05289217358955933
$ws.Range("T7").Value = 0.05289217358955934

# Row 8
$ws.Range("G8").Value = 26.10071233333333
$ws.Range("H8").Value = 78.302137
$ws.Range("I8").Value = 0.1508972550112184
$ws.Range("J8").Value = 0.1508972550112184
$ws.Range("M8").Value = 0.01339666666666667
$ws.Range("N8").Value = 0.04019
$ws.Range("O8").Value = 0.08393217762128816
$ws.Range("P8").Value = 0.08393217762128814
$ws.Range("Q8").Value = 0.3496625428922222
$ws.Range("R8").Value = 3.14696288603
$ws.Range("S8").Value = 0.0126651352101664
$ws.Range("T8").Value = 0.01266513521016639

# Row 9
$ws.Range("G9").Value = 26.10071233333333
$ws.Range("H9").Value = 78.302137
$ws.Range("I9").Value = 0.1508972550112184
$ws.Range("J9").Value = 0.1508972550112184
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.1260863333333333
$ws.Range("N9").Value = 0.378259
$ws.Range("O9").Value = 0.7899502755623498
$ws.Range("P9").Value = 0.7899502755623498
$ws.Range("Q9").Value = 3.290943115498111
$ws.Range("R9").Value = 29.618488039483
$ws.Range("S9").Value = 0.1192013281777141
$ws.Range("T9").Value = 0.1192013281777141

# Row 10
$ws.Range("G10").Value = 26.10071233333333
$ws.Range("H10").Value = 78.302137
$ws.Range("I10").Value = 0.1508972550112184
$ws.Range("J10").Value = 0.1508972550112184
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.02013
$ws.Range("N10").Value = 0.06039
$ws.Range("O10").Value = 0.1261175468163621
$ws.Range("P10").Value = 0.1261175468163621
$ws.Range("Q10").Value = 0.5254073392699999
$ws.Range("R10").Value = 4.72866605343
$ws.Range("S10").Value = 0.01903079162333786
$ws.Range("T10").Value = 0.01903079162333786

# Row 11
$ws.Range("G11").Value = 27.85999533333333
$ws.Range("H11").Value = 83.57998600000001
$ws.Range("I11").Value = 0.1610682791617304
$ws.Range("J11").Value = 0.1610682791617305
$ws.Range("M11").Value = 0.01339666666666667
$ws.Range("N11").Value = 0.04019
$ws.Range("O11").Value = 0.08393217762128816
$ws.Range("P11").Value = 0.08393217762128814
$ws.Range("Q11").Value = 0.3732310708155556
$ws.Range("R11").Value = 3.35907963734
$ws.Range("S11").Value = 0.01351881141575758
$ws.Range("T11").Value = 0.01351881141575758

# Row 12
$ws.Range("G12").Value = 27.85999533333333
$ws.Range("H12").Value = 83.57998600000001
$ws.Range("I12").Value = 0.1610682791617304
$ws.Range("J12").Value = 0.1610682791617305
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.1260863333333333
$ws.Range("N12").Value = 0.378259
$ws.Range("O12").Value = 0.7899502755623498
$ws.Range("P12").Value = 0.7899502755623498
$ws.Range("Q12").Value = 3.512764658263777
$ws.Range("R12").Value = 31.614881924374
$ws.Range("S12").Value = 0.1272359315081625
$ws.Range("T12").Value = 0.1272359315081625

# Row 13
$ws.Range("G13").Value = 27.85999533333333
$ws.Range("H13").Value = 83.57998600000001
$ws.Range("I13").Value = 0.1610682791617304
$ws.Range("J13").Value = 0.1610682791617305
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.02013
$ws.Range("N13").Value = 0.06039
$ws.Range("O13").Value = 0.1261175468163621
$ws.Range("P13").Value = 0.1261175468163621
$ws.Range("Q13").Value = 0.5608217060599999
$ws.Range("R13").Value = 5.04739535454
$ws.Range("S13").Value = 0.02031353623781041
$ws.Range("T13").Value = 0.02031353623781042

# Row 14
$ws.Range("G14").Value = 30.09393033333333
$ws.Range("H14").Value = 90.281791
$ws.Range("I14").Value = 0.1739834308659612
$ws.Range("J14").Value = 0.1739834308659612
$ws.Range("M14").Value = 0.01339666666666667
$ws.Range("N14").Value = 0.04019
$ws.Range("O14").Value = 0.08393217762128816
$ws.Range("P14").Value = 0.08393217762128814
$ws.Range("Q14").Value = 0.4031583533655556
$ws.Range("R14").Value = 3.62842518029
$ws.Range("S14").Value = 0.01460280822260296
$ws.Range("T14").Value = 0.01460280822260296

# Row 15
$ws.Range("G15").Value = 30.09393033333333
$ws.Range("H15").Value = 90.281791
$ws.Range("I15").Value = 0.1739834308659612
$ws.Range("J15").Value = 0.1739834308659612
$ws.Range("K15").Value = 1
$ws.Range("L15").Value = 0.3333333333333333
$ws.Range("M15").Value = 0.1260863333333333
$ws.Range("N15").Value = 0.378259
$ws.Range("O15").Value = 0.7899502755623498
$ws.Range("P15").Value = 0.7899502755623498
$ws.Range("Q15").Value = 3.794433331318777
$ws.Range("R15").Value = 34.149899981869
$ws.Range("S15").Value = 0.1374382591558491
$ws.Range("T15").Value = 0.1374382591558491

# Row 16
$ws.Range("G16").Value = 30.09393033333333
$ws.Range("H16").Value = 90.281791
$ws.Range("I16").Value = 0.1739834308659612
$ws.Range("J16").Value = 0.1739834308659612
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.02013
$ws.Range("N16").Value = 0.06039
$ws.Range("O16").Value = 0.1261175468163621
$ws.Range("P16").Value = 0.1261175468163621
$ws.Range("Q16").Value = 0.60579081761
$ws.Range("R16").Value = 5.45211735849
$ws.Range("S16").Value = 0.02194236348750915
$ws.Range("T16").Value = 0.02194236348750915
